# Mark the e716d8a8-... entry as "Ready for handoff" (a new handoff report
# was generated for it), updating the Overview roll-up sheet and the
# per-locale (zh-cn / de-de) detail sheets with the new handoff timestamp.

$wb = $excel.ActiveWorkbook

# Overview sheet: Status column mirrors into both the zh-cn (col B) and
# de-de (col C) summary columns for this row.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# zh-cn detail sheet: Status + Latest Handoff Datetime for the same file.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "2016-01-13 11:36:06"

# de-de detail sheet: Status + Latest Handoff Datetime for the same file.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "2016-01-13 11:36:30"
